$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column C header from "Tipo" to "Tipo de prueba"
$ws.Range("C1").Value = "Tipo de prueba"

# Delete the "Descripción del Caso de Prueba" column (column E) entirely,
# shifting all columns to its right one position to the left
$ws.Columns("E").Delete()

# Give the "Tipo de prueba" column an explicit width (closest achievable
# value given the engine's internal pixel-width quantization)
$ws.Columns("C").ColumnWidth = 16.65

# Reset the active selection to mirror the upstream edit
$ws.Range("E1:E1048576").Select() | Out-Null
